$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.4173876666666667
$ws.Range("H2").Value = 1.252163
$ws.Range("I2").Value = 0.1865415014963835
$ws.Range("J2").Value = 0.1865415014963835
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.06449866666666666
$ws.Range("N2").Value = 0.193496
$ws.Range("O2").Value = 0.004525829983623641
$ws.Range("P2").Value = 0.004525829983623642
$ws.Range("Q2").Value = 0.02692094798311111
$ws.Range("R2").Value = 0.242288531848
$ws.Range("S2").Value = 0.0008442551206625069
$ws.Range("T2").Value = 0.0008442551206625071

# Row 3
$ws.Range("G3").Value = 0.4173876666666667
$ws.Range("H3").Value = 1.252163
$ws.Range("I3").Value = 0.1865415014963835
$ws.Range("J3").Value = 0.1865415014963835
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 10.61985133333333
$ws.Range("N3").Value = 31.859554
$ws.Range("O3").Value = 0.745188142173877
$ws.Range("P3").Value = 0.7451881421738772
$ws.Range("Q3").Value = 4.43259496836689
$ws.Range("R3").Value = 39.893354715302
$ws.Range("S3").Value = 0.1390085149384155
$ws.Range("T3").Value = 0.1390085149384156

# Row 4
$ws.Range("G4").Value = 0.4173876666666667
$ws.Range("H4").Value = 1.252163
$ws.Range("I4").Value = 0.1865415014963835
$ws.Range("J4").Value = 0.1865415014963835
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.566885000000001
$ws.Range("N4").Value = 10.700655
$ws.Range("O4").Value = 0.2502860278424993
$ws.Range("P4").Value = 0.2502860278424993
$ws.Range("Q4").Value = 1.488773807418334
$ws.Range("R4").Value = 13.398964266765
$ws.Range("S4").Value = 0.04668873143730547
$ws.Range("T4").Value = 0.04668873143730548

# Row 5
$ws.Range("G5").Value = 1.220831
$ws.Range("H5").Value = 3.662493
$ws.Range("I5").Value = 0.5456214114616024
$ws.Range("J5").Value = 0.5456214114616023
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.06449866666666666
$ws.Range("N5").Value = 0.193496
$ws.Range("O5").Value = 0.004525829983623641
$ws.Range("P5").Value = 0.004525829983623642
$ws.Range("Q5").Value = 0.07874197172533334
$ws.Range("R5").Value = 0.7086777455280001
$ws.Range("S5").Value = 0.002469389743699972
$ws.Range("T5").Value = 0.002469389743699972

# Row 6
$ws.Range("G6").Value = 1.220831
$ws.Range("H6").Value = 3.662493
$ws.Range("I6").Value = 0.5456214114616024
$ws.Range("J6").Value = 0.5456214114616023
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 10.61985133333333
$ws.Range("N6").Value = 31.859554
$ws.Range("O6").Value = 0.745188142173877
$ws.Range("P6").Value = 0.7451881421738772
$ws.Range("Q6").Value = 12.96504372312467
$ws.Range("R6").Value = 116.685393508122
$ws.Range("S6").Value = 0.40659060593736
$ws.Range("T6").Value = 0.40659060593736

# Row 7
$ws.Range("G7").Value = 1.220831
$ws.Range("H7").Value = 3.662493
$ws.Range("I7").Value = 0.5456214114616024
$ws.Range("J7").Value = 0.5456214114616023
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.566885000000001
$ws.Range("N7").Value = 10.700655
$ws.Range("O7").Value = 0.2502860278424993
$ws.Range("P7").Value = 0.2502860278424993
$ws.Range("Q7").Value = 4.354563781435002
$ws.Range("R7").Value = 39.19107403291501
$ws.Range("S7").Value = 0.1365614157805423
$ws.Range("T7").Value = 0.1365614157805423

# Row 8
$ws.Range("G8").Value = 0.599287
$ws.Range("H8").Value = 1.797861
$ws.Range("I8").Value = 0.2678370870420142
$ws.Range("J8").Value = 0.2678370870420142
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.06449866666666666
$ws.Range("N8").Value = 0.193496
$ws.Range("O8").Value = 0.004525829983623641
$ws.Range("P8").Value = 0.004525829983623642
$ws.Range("Q8").Value = 0.03865321245066666
$ws.Range("R8").Value = 0.347878912056
$ws.Range("S8").Value = 0.001212185119261163
$ws.Range("T8").Value = 0.001212185119261163

# Row 9
$ws.Range("G9").Value = 0.599287
$ws.Range("H9").Value = 1.797861
$ws.Range("I9").Value = 0.2678370870420142
$ws.Range("J9").Value = 0.2678370870420142
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 10.61985133333333
$ws.Range("N9").Value = 31.859554
$ws.Range("O9").Value = 0.745188142173877
$ws.Range("P9").Value = 0.7451881421738772
$ws.Range("Q9").Value = 6.364338845999334
$ws.Range("R9").Value = 57.279049613994
$ws.Range("S9").Value = 0.1995890212981015
$ws.Range("T9").Value = 0.1995890212981016

# Row 10
$ws.Range("G10").Value = 0.599287
$ws.Range("H10").Value = 1.797861
$ws.Range("I10").Value = 0.2678370870420142
$ws.Range("J10").Value = 0.2678370870420142
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.566885000000001
$ws.Range("N10").Value = 10.700655
$ws.Range("O10").Value = 0.2502860278424993
$ws.Range("P10").Value = 0.2502860278424993
$ws.Range("Q10").Value = 2.137587810995
$ws.Range("R10").Value = 19.238290298955
$ws.Range("S10").Value = 0.06703588062465146
$ws.Range("T10").Value = 0.06703588062465148
